$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.859.22"
$ws.Range("E2").Value = "  +0.61%  "
$ws.Range("D3").Value = "2.679.80"
$ws.Range("E3").Value = "  +0.71%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'604.00"
$ws.Range("E5").Value = "  -0.10%  "
$ws.Range("D6").Value = "'156.57"
$ws.Range("E6").Value = "  -0.75%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.603"
$ws.Range("E8").Value = "  +2.39%  "
$ws.Range("D9").Value = "'0.122"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").Value = "'5.95"
$ws.Range("E10").Value = "  +2.40%  "
$ws.Range("D11").Value = "'0.398"
$ws.Range("E11").Value = "  -2.43%  "
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("D13").Value = "'29.46"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "'0.0000197"
$ws.Range("E14").Value = "  +3.00%  "
$ws.Range("D15").Value = "3.161.49"
$ws.Range("E15").Value = "  +0.67%  "
$ws.Range("D16").Value = "65.669.43"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").Value = "2.701.49"
$ws.Range("E17").Value = "  +0.67%  "
$ws.Range("D18").Value = "'12.59"
$ws.Range("E18").Value = "  -1.49%  "
$ws.Range("D19").Value = "'4.83"
$ws.Range("E19").Value = "  -1.28%  "
$ws.Range("D20").Value = "'7.54"
$ws.Range("E20").Value = "  +2.35%  "
$ws.Range("D21").Value = "'351.42"
$ws.Range("E21").Value = "  -2.39%  "
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").Value = "'70.55"
$ws.Range("E23").Value = "  +2.38%  "
$ws.Range("E24").Value = "  +6.51%  "
$ws.Range("D25").Value = "'9.85"
$ws.Range("E25").Value = "  +3.50%  "
$ws.Range("D26").Value = "'1.63"
$ws.Range("E26").Value = "  -4.37%  "
$ws.Range("E27").Value = "  -1.61%  "
$ws.Range("E28").Value = "  +1.85%  "
$ws.Range("D29").Value = "'8.15"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("E31").Value = "  -1.56%  "
$ws.Range("D32").Value = "'535.14"
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").Value = "'1.77"
$ws.Range("E33").Value = "  -3.92%  "
$ws.Range("D34").Value = "'6.56"
$ws.Range("E34").Value = "  +2.13%  "
$ws.Range("D35").Value = "'5.39"
$ws.Range("E35").Value = "  -4.57%  "
$ws.Range("D36").Value = "'0.425"
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("D37").Value = "'20.43"
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").Value = "'160.41"
$ws.Range("E38").Value = "  -1.57%  "
$ws.Range("D40").Value = "'1.96"
$ws.Range("E40").Value = "  -2.65%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  -0.52%  "
$ws.Range("D43").Value = "'166.00"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'4.10"
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "'0.0618"
$ws.Range("E45").Value = "  +0.32%  "
$ws.Range("D46").Value = "'23.13"
$ws.Range("E46").Value = "  -0.17%  "
$ws.Range("D47").Value = "'2.23"
$ws.Range("E47").Value = "  -3.85%  "
$ws.Range("D48").Value = "'0.0263"
$ws.Range("E48").Value = "  -0.62%  "
$ws.Range("D49").Value = "'0.651"
$ws.Range("E49").Value = "  -1.82%  "
$ws.Range("D50").Value = "'20.26"
$ws.Range("E50").Value = "  +2.49%  "
$ws.Range("D51").Value = "'0.1000"
$ws.Range("E51").Value = "  +1.30%  "
